# Applies the Tue Apr 16 09:57:30 UTC 2024 "Updated cryptos list" data refresh.
# Price (column D) and Volume(1h) (column E) values are plain text cells in the
# source sheet (numbers rendered with "." thousands separators, fixed decimal
# widths, etc.), so numeric-looking replacements are written with a leading
# apostrophe to force Excel to store them as text (matching the original
# inlineStr cell type) instead of silently re-parsing/reformatting them as
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.278.54'
$ws.Range("E2").Value = '  -5.15%  '

$ws.Range("D3").Value = '''3.078.97'
$ws.Range("E3").Value = '  -5.52%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = '''546.22'
$ws.Range("E5").Value = '  -6.32%  '

$ws.Range("D6").Value = '''135.80'
$ws.Range("E6").Value = '  -12.12%  '

$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("D8").Value = '''3.071.51'
$ws.Range("E8").Value = '  -5.46%  '

$ws.Range("D9").Value = '''0.493'
$ws.Range("E9").Value = '  -4.34%  '

$ws.Range("D10").Value = '''0.155'
$ws.Range("E10").Value = '  -6.05%  '

$ws.Range("D11").Value = '''6.21'
$ws.Range("E11").Value = '  -12.45%  '

$ws.Range("D12").Value = '''0.467'
$ws.Range("E12").Value = '  -4.64%  '

$ws.Range("D13").Value = '''35.11'
$ws.Range("E13").Value = '  -7.58%  '

$ws.Range("E14").Value = '  -8.54%  '

$ws.Range("D15").Value = '''3.583.55'
$ws.Range("E15").Value = '  -5.38%  '

$ws.Range("D16").Value = '''63.315.79'
$ws.Range("E16").Value = '  -5.18%  '

$ws.Range("E17").Value = '  -3.38%  '

$ws.Range("D18").Value = '''3.084.69'
$ws.Range("E18").Value = '  -5.31%  '

$ws.Range("D19").Value = '''6.72'
$ws.Range("E19").Value = '  -6.09%  '

$ws.Range("D20").Value = '''485.29'
$ws.Range("E20").Value = '  -13.23%  '

$ws.Range("D21").Value = '''13.52'
$ws.Range("E21").Value = '  -6.77%  '

$ws.Range("D22").Value = '''0.714'
$ws.Range("E22").Value = '  -4.34%  '

$ws.Range("D23").Value = '''7.23'
$ws.Range("E23").Value = '  -7.36%  '

$ws.Range("D24").Value = '''78.73'
$ws.Range("E24").Value = '  -4.04%  '

$ws.Range("D25").Value = '''12.28'
$ws.Range("E25").Value = '  -10.38%  '

$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("D27").Value = '''8.45'
$ws.Range("E27").Value = '  -8.94%  '

$ws.Range("D28").Value = '''2.74'
$ws.Range("E28").Value = '  -8.06%  '

$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("D30").Value = '''1.95'
$ws.Range("E30").Value = '  -13.00%  '

$ws.Range("D31").Value = '''26.50'
$ws.Range("E31").Value = '  -5.15%  '

$ws.Range("E32").Value = '  -4.96%  '

$ws.Range("D33").Value = '''2.49'
$ws.Range("E33").Value = '  -10.08%  '

$ws.Range("D34").Value = '''59.36'
$ws.Range("E34").Value = '  +7.11%  '

$ws.Range("D35").Value = '''500.27'
$ws.Range("E35").Value = '  -11.02%  '

$ws.Range("D36").Value = '''6.02'
$ws.Range("E36").Value = '  -6.01%  '

$ws.Range("D37").Value = '''5.07'

$ws.Range("D38").Value = '''3.144.37'
$ws.Range("E38").Value = '  -1.38%  '

$ws.Range("D39").Value = '''0.0395'
$ws.Range("E39").Value = '  -13.85%  '

$ws.Range("D40").Value = '''0.0797'
$ws.Range("E40").Value = '  -8.05%  '

$ws.Range("D41").Value = '''0.118'
$ws.Range("E41").Value = '  -10.61%  '

$ws.Range("D42").Value = '''8.12'
$ws.Range("E42").Value = '  -6.17%  '

$ws.Range("D43").Value = '''2.59'
$ws.Range("E43").Value = '  -15.67%  '

$ws.Range("D44").Value = '''0.254'
$ws.Range("E44").Value = '  -9.04%  '

$ws.Range("D45").Value = '''0.999'
$ws.Range("E45").Value = '  +0.02%  '

$ws.Range("D46").Value = '''25.23'
$ws.Range("E46").Value = '  -4.70%  '

$ws.Range("D47").Value = '''2.05'
$ws.Range("E47").Value = '  -12.09%  '

$ws.Range("D48").Value = '''119.81'
$ws.Range("E48").Value = '  -5.16%  '

$ws.Range("D49").Value = '''0.108'
$ws.Range("E49").Value = '  -4.45%  '

$ws.Range("D50").Value = '''0.0₃0506'
$ws.Range("E50").Value = '  -9.62%  '

$ws.Range("B51").Value = 'CoreDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D51").Value = '''2.31'
$ws.Range("E51").Value = '  +29.19%  '
